# Atualização de bases das ligas, do dia: 14-06-2024 às 20:31
#
# The source feed re-sorted a handful of fixtures that share the same
# Date/Div, so the two data rows for each of these pairs need to swap
# places. The leading row-index column (A) stays put; every other
# column (B..AD: id, HomeTeam, AwayTeam, scores, odds, etc.) swaps
# between the two rows of each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # column B
$lastCol  = 30  # column AD

$rowPairs = @(
    @(32, 33),
    @(37, 38),
    @(67, 68),
    @(246, 247),
    @(252, 253)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $val1 = $cell1.Value()
        $val2 = $cell2.Value()

        $cell1.Value = $val2
        $cell2.Value = $val1
    }
}
